# mau_import_loai_thiet_bi.xlsx
# cap nhat importexcel taisan/vitri, fix dm loai bao tri, ke hoach bao tri
#
# The "Trực thuộc" (parent link) column (E) for the two sub-type rows
# (Máy tính / Máy in) used to hold stray standalone codes "VT03"/"VT02"
# that didn't correspond to anything else in the sheet. They are fixed
# up to point at the real parent "Mã Loại" codes from column B
# (LTB001 = Xe tải, LTB002 = Xe container). The leftover placeholder
# single-space value in D5 is also cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# D5 previously held a lone space " " used as a placeholder - clear it.
$ws.Range("D5").Value = ""

# E6/E7: replace the orphan "VT03"/"VT02" codes with the correct parent
# "Mã Loại" references (LTB001 = Xe tải, LTB002 = Xe container).
$ws.Range("E6").Value = "LTB001"
$ws.Range("E7").Value = "LTB002"

# Give them the same centered/bordered look as the neighbouring "x"
# cells in column D, instead of the old un-centered style.
$ws.Range("E6:E7").HorizontalAlignment = $xlCenter
$ws.Range("E6:E7").VerticalAlignment = $xlCenter
